# Staff Details workbook - rework the "Assigned Streams" column into a set
# of per-category "stream" columns (Production, Prescription Scanning,
# Legacy, Portering, Student Services, Scanning Services), populate every
# staff member's streams under the right category, and drop the old
# "legacy stream controls" single-column data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 2: replace "Assigned Streams" (AD2) with the new category
#     headers spread across AD2:AI2 ---------------------------------------
$ws.Range("AD2").Value = "Production"
$ws.Range("AE2").Value = "Prescription Scanning"
$ws.Range("AF2").Value = "Legacy"
$ws.Range("AG2").Value = "Portering"
$ws.Range("AH2").Value = "Student Services"
$ws.Range("AI2").Value = "Scanning Services"

# --- Row 3 (Paul Smith) ---------------------------------------------------
$ws.Range("AD3").Value = "Islands Account, Islands Reconciliation"
$ws.Range("AE3").ClearContents()
$ws.Range("AF3").Value = "BTST"
$ws.Range("AG3").Value = "Bailing"
$ws.Range("AH3").ClearContents()
$ws.Range("AI3").ClearContents()

# --- Row 4 (Paul Jones) ---------------------------------------------------
$ws.Range("AD4").Value = "Item Confirmation"
$ws.Range("AE4").ClearContents()
$ws.Range("AF4").ClearContents()
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").Value = "Scan Applications, Scan CCR"
$ws.Range("AI4").ClearContents()

# --- Row 5 (Tony Robinson) ------------------------------------------------
$ws.Range("AD5").ClearContents()
$ws.Range("AE5").Value = "PADM - Docketing, PADM - Scanning"
$ws.Range("AF5").Value = "Applications"
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AI5").ClearContents()

# --- Column widths: the old single wide "AD" column (82.71) is replaced by
#     six narrower equal-width columns AD:AI (~40.71) ----------------------
$ws.Range("AD1:AI1").ColumnWidth = 39.8

# --- AutoFilter / defined name now need to span out to column AI ---------
$ws.AutoFilterMode = $false
$ws.Range("A2:AI5").AutoFilter()
$wb.Names.Item(1).RefersTo = "='Staff Details'!`$A`$2:`$AI`$5"
